$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates
$ws.Range("M1").Value = "CADD_phred"
$ws.Range("P1").Value = "cosmic86_coding"

# VHL row (row 2) updates
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "50.98%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "52|102"
